$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 169: 2020-02-27
$ws.Cells.Item(169, 1).Value = 1582761600
Set-TextCell 169 2 "2020-02-27"
Set-TextCell 169 3 "0208"
$ws.Cells.Item(169, 4).Value = "GREATEC"
$ws.Cells.Item(169, 5).Value = 3.68
$ws.Cells.Item(169, 6).Value = 3.68
$ws.Cells.Item(169, 7).Value = 3.35
$ws.Cells.Item(169, 8).Value = 3.41
$ws.Cells.Item(169, 9).Value = 5779700

# Row 170: 2020-02-28
$ws.Cells.Item(170, 1).Value = 1582848000
Set-TextCell 170 2 "2020-02-28"
Set-TextCell 170 3 "0208"
$ws.Cells.Item(170, 4).Value = "GREATEC"
$ws.Cells.Item(170, 5).Value = 3.3
$ws.Cells.Item(170, 6).Value = 3.35
$ws.Cells.Item(170, 7).Value = 3.06
$ws.Cells.Item(170, 8).Value = 3.23
$ws.Cells.Item(170, 9).Value = 8964700
